$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.355.03"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.34%  '

$ws.Range("D3").Value = "'2.177.22"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.47%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").Value = "'252.54"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.09%  '

$ws.Range("D6").Value = "'0.610"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.53%  '

$ws.Range("D7").Value = "'74.42"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.67%  '

$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("D9").Value = "'0.578"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -5.43%  '

$ws.Range("D10").Value = "'39.96"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.35%  '

$ws.Range("D11").Value = "'0.0909"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.58%  '

$ws.Range("E12").Value = '  -0.80%  '

$ws.Range("D13").Value = "'6.71"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.61%  '

$ws.Range("D14").Value = "'2.504.29"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.37%  '

$ws.Range("D15").Value = "'14.14"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.37%  '

$ws.Range("D16").Value = "'2.171.31"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.24%  '

$ws.Range("D17").Value = "'0.766"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -6.70%  '

$ws.Range("D18").Value = "'42.280.55"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.19%  '

$ws.Range("D19").Value = "'0.0000101"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.88%  '

$ws.Range("D20").Value = "'70.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.32%  '

$ws.Range("D21").Value = "'5.84"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.25%  '

$ws.Range("D22").Value = "'225.20"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.42%  '

$ws.Range("D23").Value = "'9.37"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -12.22%  '

$ws.Range("D24").Value = "'2.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.00%  '

$ws.Range("E25").Value = '  +0.00%  '

$ws.Range("D26").Value = "'10.39"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -5.38%  '

$ws.Range("D27").Value = "'3.36"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.18%  '

$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").Value = "'2.18"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.36%  '

$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D29").Value = "'2.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.67%  '

$ws.Range("D30").Value = "'37.53"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.13%  '

$ws.Range("D31").Value = "'171.46"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.45%  '

$ws.Range("D32").Value = "'19.98"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.90%  '

$ws.Range("D33").Value = "'0.0825"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.46%  '

$ws.Range("D34").Value = "'5.13"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.90%  '

$ws.Range("E35").Value = '  -2.13%  '

$ws.Range("D36").Value = "'0.107"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.90%  '

$ws.Range("D37").Value = "'0.0336"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.19%  '

$ws.Range("D38").Value = "'4.16"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.11%  '

$ws.Range("D39").Value = "'11.87"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -10.93%  '

$ws.Range("D40").Value = "'2.05"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.52%  '

$ws.Range("E41").Value = '  +12.64%  '

$ws.Range("D42").Value = "'0.194"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.20%  '

$ws.Range("B43").Value = 'THORChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D43").Value = "'5.14"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -7.51%  '

$ws.Range("B44").Value = 'MultiversX'
$ws.Range("C44").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D44").Value = "'58.58"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.22%  '

$ws.Range("D45").Value = "'101.12"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.17%  '

$ws.Range("D46").Value = "'0.0969"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.54%  '

$ws.Range("B47").Value = 'FraxShare'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D47").Value = "'8.13"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.50%  '

$ws.Range("B48").Value = 'WOONetwork'
$ws.Range("C48").Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range("D48").Value = "'0.455"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.12%  '

$ws.Range("D49").Value = "'1.08"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.99%  '

$ws.Range("D50").Value = "'1.12"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.26%  '

$ws.Range("E51").Value = '  +0.17%  '
